$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5 for "llama3_70b_instruct" (keeps the list
# alphabetically/logically ordered) and pushes all following rows down by one.
$ws.Rows.Item(5).Insert()

$ws.Cells.Item(5, 1).Value = "llama3_70b_instruct"
$ws.Cells.Item(5, 2).Value = 84.26
$ws.Cells.Item(5, 3).Value = 71.23
$ws.Cells.Item(5, 4).Value = 71.89
$ws.Cells.Item(5, 5).Value = 66.86
$ws.Cells.Item(5, 6).Value = 9.98
$ws.Cells.Item(5, 7).Value = 8.15
$ws.Cells.Item(5, 8).Value = 6.79
$ws.Cells.Item(5, 9).Value = 0.04

# Round all pre-existing percentage metric values (columns B:I) down to 2
# decimal places to make the metrics file robust / consistent with the
# newly added JSON-derived details.
$data = @{
    2  = @(35.16, 5.91, 7.52, 14.97, 12.33, 30.31, 39.38, 35.16)
    3  = @(67.74, 45.76, 55.82, 50.86, 50.64, 29.69, 51.01, 26.72)
    4  = @(56.11, 24.77, 24.62, 19.74, 28.48, 22.39, 30.79, 23.16)
    6  = @(53.72, 48.77, 34.75, 31.23, 2.35, 3.16, 6.86, 13.5)
    7  = @(56.88, 32.33, 25.28, 30.42, 41.03, 30.28, 36, 36.4)
    8  = @(83.34, 43.34, 52.29, 44.37, 64.33, 75.27, 54.86, 66.2)
    9  = @(91.23, 38.86, 47.52, 39.27, 5.69, 3.93, 2.42, 2.64)
    10 = @(89.28, 55.08, 55.08, 47.63, 26.46, 53.21, 38.06, 54.02)
    11 = @(86.31, 68.33, 59.23, 63.82, 0.37, 1.98, 0.99, 3.23)
    12 = @(94.39, 66.79, 67.3, 65.32, 0.99, 6.2, 2.09, 9.21)
    13 = @(49.98, 46.68, 38.86, 43.63, 45.91, 12.95, 15.6, 9.06)
    14 = @(35.45, 30.5, 21.69, 26.97, 19.23, 19.49, 43.19, 36.7)
    15 = @(66.79, 29.83, 20.15, 21.69, 47.6, 44.04, 43.41, 36.51)
    16 = @(79.08, 42.02, 33.76, 47.41, 68.99, 59.41, 59.27, 44.62)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($r, $i + 2).Value = $vals[$i]
    }
}
